# Saldo.xlsx update ("Add files via upload"):
#  - Remove the "004329030 / DANIELA / 17759.85" row entirely.
#  - Move the 8-row block "004468717 HELOISA" .. "004452597 LARA"
#    (which sat right after "004487016 ROGERIO") to sit right BEFORE
#    the ROGERIO row instead.
#  - Change ROGERIO's balance from 2982.44 to 870.94 (now the last row
#    of that block).
#
# The two account/name columns are plain text (leading zeros must be
# kept), so every cross-cell copy below uses Range.Copy(destination)
# rather than re-typing values through .Value/.Value2 (which would
# re-interpret "004468717" etc. as a number and strip the leading
# zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the anchor rows by account number instead of hard-coded row
# indices, so the script is resilient to the sheet's exact layout.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

$danielaRow = 0
$rogerioRow = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $acct = $ws.Cells.Item($r, 1).Value2
    if ($acct -eq "004329030") { $danielaRow = $r }
    if ($acct -eq "004487016") { $rogerioRow = $r }
}

if ($danielaRow -eq 0 -or $rogerioRow -eq 0) {
    throw "Could not locate anchor rows (DANIELA=$danielaRow, ROGERIO=$rogerioRow)"
}

# 1) Delete DANIELA's row outright.
$ws.Rows.Item($danielaRow).Delete()

# Deleting a row above ROGERIO shifts it up by one.
if ($danielaRow -lt $rogerioRow) {
    $rogerioRow = $rogerioRow - 1
}

# The block to relocate is the 8 rows directly under ROGERIO.
$blockFirst = $rogerioRow + 1
$blockLast = $rogerioRow + 8
$blockHeight = $blockLast - $blockFirst + 1

# Scratch area, far away from any real data, used as a holding pen so
# the ROGERIO row and the block don't clobber each other while moving.
$scratchBlock = $ws.Range("A2000:C2007")
$scratchRogerio = $ws.Range("A2008:C2008")

# 2) Stash the HELOISA..LARA block, then the ROGERIO row.
$ws.Range($ws.Cells.Item($blockFirst, 1), $ws.Cells.Item($blockLast, 3)).Copy($scratchBlock)
$ws.Range($ws.Cells.Item($rogerioRow, 1), $ws.Cells.Item($rogerioRow, 3)).Copy($scratchRogerio)

# 3) Write the block back one row earlier (into ROGERIO's old slot and
#    onward), then put ROGERIO back after it.
$destBlock = $ws.Range($ws.Cells.Item($rogerioRow, 1), $ws.Cells.Item($rogerioRow + $blockHeight - 1, 3))
$scratchBlock.Copy($destBlock)

$newRogerioRow = $rogerioRow + $blockHeight
$destRogerio = $ws.Range($ws.Cells.Item($newRogerioRow, 1), $ws.Cells.Item($newRogerioRow, 3))
$scratchRogerio.Copy($destRogerio)

# 4) Update ROGERIO's balance at its new location.
$ws.Cells.Item($newRogerioRow, 3).Value2 = 870.94

# 5) Clean up the scratch area so it doesn't affect the used range.
$ws.Range("A2000:C2008").Clear()
